$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.049.15'
$ws.Range('E2').Value = '  -2.51%  '
$ws.Range('D3').Value = '1.821.97'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('E4').Value = '  -0.93%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.15'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4224'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.86%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3666'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07217'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8405'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.26%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.74'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.93%  '
$ws.Range('D12').Value = '1.827.81'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.669'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.07073'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.287'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.66%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '90.19'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008747'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.69%  '
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.91'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.57%  '
$ws.Range('D21').Value = '27.124.19'
$ws.Range('E21').Value = '  -2.24%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.138'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.85'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.97%  '
$ws.Range('D24').Value = '2.048.51'
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.979'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '152.04'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('E27').Value = '  +3.26%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.24'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.277'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '117.29'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08737'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.23%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.177'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.30%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7369'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.32%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.906'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.64%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.419'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.82%  '
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.088'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.79%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01951'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05247'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '7.322'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.868'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.58%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1689'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.5040'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.554'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.93%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.54'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.86%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '106.16'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4706'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.64%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.922'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.30%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.9999'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06334'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.648'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.39%  '
